$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the multi-line JSON text that goes into the new G9 cell / new shared string.
$lines = @(
  'Download',
  '{',
  '  "extracted_name": "Butane",',
  '  "matched_name": "Butane",',
  '  "matching_score": 1,',
  '  "reason": "Exact name match",',
  '  "etrm_code": "400007248.0",',
  '  "etrm_id": "1560",',
  '  "alternatives": [',
  '    {',
  '      "matched_name": "BUTANE/BUTLYENE SALES",',
  '      "score": 0.8,',
  '      "code": "400007558.0",',
  '      "id": "2251"',
  '    },',
  '    {',
  '      "matched_name": "Normal Butane Frac",',
  '      "score": 0.8,',
  '      "code": "400007561.0",',
  '      "id": "2180"',
  '    },',
  '    {',
  '      "matched_name": "Refridg Normal Butane",',
  '      "score": 0.8,',
  '      "code": "400007562.0",',
  '      "id": "2071"',
  '    },',
  '    {',
  '      "matched_name": "Cap and Trade Butane Obligations Quebec",',
  '      "score": 0.8,',
  '      "code": null,',
  '      "id": "2691"',
  '    },',
  '    {',
  '      "matched_name": "Cap and Trade Butane Potential Obligations Quebec",',
  '      "score": 0.8,',
  '      "code": null,',
  '      "id": "2692"',
  '    }',
  '  ]',
  '}'
)
$newText = [string]::Join("`n", $lines)

# New column G holds the extra matching details for the Butane row (row 9).
$ws.Range("G9").Value = $newText
$ws.Range("G9").WrapText = $true

# Give column G the same generous width used for the rest of the long text columns
# (~82.44 characters wide).
$ws.Columns.Item(7).ColumnWidth = 81.65

# Row 9 grows slightly less tall than the other (overflowing) rows to fit the new content.
$ws.Rows.Item(9).RowHeight = 408.6

# Reset the view: scroll back to the top and select F1 (matches the saved view state).
$ws.Range("F1").Select()
